$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Scum Master
$ws.Range("A2").Value = "marcollano5@hotmail.com"
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:marcollano5@hotmail.com")
$ws.Range("B2").Value = "Scum Master"

# Row 3: Product Owner
$ws.Range("A3").Value = "marco_llano5@hotmail.com"
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:marco_llano5@hotmail.com")
$ws.Range("B3").Value = "Product Owner"

$ws.Range("A4").Select() | Out-Null

Write-Output "done"
